$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the long step-by-step instruction texts in column D (rows 2-5)
# with short action titles, per the updated test-scenario template.
$ws.Range("D2").Value = "Tambah Setup Profile Bank"
$ws.Range("D3").Value = "View Setup Profile Bank"
$ws.Range("D4").Value = "Ubah Setup Profile Bank"
$ws.Range("D5").Value = "Hapus Setup Profile Bank"

# Row heights shrink now that the cells hold one short line instead of
# several numbered steps.
$ws.Rows.Item(2).RowHeight = 30
$ws.Rows.Item(3).RowHeight = 30
$ws.Rows.Item(4).EntireRow.AutoFit()
$ws.Rows.Item(5).RowHeight = 30
